# Adding data & run-mode in excel.
#
# - TestCases!A4 test-case name is renamed from "checkDuplicateUPN" to
#   "checkUPN" (matching the "checkUPN" worksheet it drives).
# - The "TestCases" sheet/tab becomes the active (selected) tab again,
#   with cell B6 selected, instead of the "checkUPN" sheet with B10
#   selected.

$wb = $excel.ActiveWorkbook

$testCases = $wb.Worksheets.Item("TestCases")

# Rename the test case referenced in row 4 of TestCases.
$testCases.Range("A4").Value = "checkUPN"

# Make "TestCases" the active sheet/tab again, with B6 selected
# (previously "checkUPN" was the active tab, with B10 selected on
# TestCases left over from a prior selection).
[void]$testCases.Activate()
[void]$testCases.Range("B6").Select()
